$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D are price strings that sometimes look like plain numbers
# (e.g. "556.38"). Force those specific cells to Text format first so Excel
# stores the exact original string instead of converting it to a floating
# point number. Cells whose new text is unambiguously non-numeric (contains
# multiple "." thousand separators, subscript digits, etc.) are left as-is,
# matching how the source workbook stores them (plain text, default style).
$ws.Range("D2").Value = "63.911.59"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").Value = "3.053.00"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.38"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.33"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.049.66"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("E11").Value = "  -9.92%  "
$ws.Range("E12").Value = "  +8.10%  "
$ws.Range("E13").Value = "  +5.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.14"
$ws.Range("E14").Value = "  +4.42%  "
$ws.Range("D15").Value = "3.546.84"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "63.952.78"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "3.052.66"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.24"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("E23").Value = "  +6.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.25"
$ws.Range("E24").Value = "  +14.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.02"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.22"
$ws.Range("E31").Value = "  +4.00%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.42"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.73"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "442.49"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0805"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +13.71%  "
$ws.Range("D41").Value = "2.979.17"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.62"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("E45").Value = "  +6.16%  "
$ws.Range("E46").Value = "  +9.23%  "
$ws.Range("E48").Value = "  +4.85%  "
$ws.Range("D49").Value = "0.0₃0515"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.27"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +3.59%  "
